# Update Kehamda quarterly income-statement workbook:
#  - drop the oldest quarter column, shift all quarters one column to the
#    left, and append the newly-published quarter (فصل سوم منتهی به 1401/12)
#    in column M, together with its release-date label and figures.
#  - the "Capital" row's two newest cells (L26:M26) lose their special
#    comma-formatted style and fall back to the plain style used by the
#    rest of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8: quarter/period header labels (columns D..M)
# ---------------------------------------------------------------------
$periodLabels = @(
    "فصل دوم منتهی به 1399/09",
    "فصل سوم منتهی به 1399/12",
    "فصل چهارم منتهی به 1400/03",
    "فصل اول منتهی به 1400/06",
    "فصل دوم منتهی به 1400/09",
    "فصل سوم منتهی به 1400/12",
    "فصل چهارم منتهی به 1401/03",
    "فصل اول منتهی به 1401/06",
    "فصل دوم منتهی به 1401/09",
    "فصل سوم منتهی به 1401/12"
)
for ($i = 0; $i -lt $periodLabels.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $periodLabels[$i]
}

# ---------------------------------------------------------------------
# Row 9: release-date labels (columns D..M)
# ---------------------------------------------------------------------
$dateLabels = @(
    "1400-11-27 (4)",
    "1401-01-30 (2)",
    "1401-07-20 (9)",
    "1401-07-30 (2)",
    "1401-12-03 (4)",
    "1402-01-30 (2)",
    "1402-01-30 (7)",
    "1401-07-30",
    "1401-12-03 (2)",
    "1402-01-30"
)
for ($i = 0; $i -lt $dateLabels.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $dateLabels[$i]
}

# ---------------------------------------------------------------------
# Data rows: shift existing quarterly figures one column left and add the
# newly-reported quarter's figures in column M.
# ---------------------------------------------------------------------
$rowData = @{
    11 = @(935948,1155829,1385903,1700675,1637371,1916177,1903047,2342887,2112863,3755932)
    12 = @(-394858,-508432,-531891,-631848,-594064,-754538,-812192,-1060461,-1168958,-2114339)
    13 = @(541090,647397,854012,1068827,1043307,1161639,1090855,1282426,943905,1641593)
    14 = @(-44172,-34108,-31128,-34938,-68397,-53129,-65347,-85260,-121067,-124770)
    16 = @(32343,-10926,-7145,0,196680,-175372,-718,0,230143,-172663)
    17 = @(529261,602363,815739,1033889,1171590,933138,1024790,1197166,1052981,1344160)
    18 = @(-2319,-3123,-3581,-55985,-56512,-57823,-78272,-173075,-237381,-368726)
    19 = @(574369,-189016,-42056,65801,124249,158720,105511,96027,630518,483212)
    20 = @(1101311,410224,770102,1043705,1239327,1034035,1052029,1120118,1446118,1458646)
    21 = @(-192329,-79643,24003,-189385,-111939,-147895,-189948,-195265,77861,-275865)
    22 = @(908982,330581,794105,854320,1127388,886140,862081,924853,1523979,1182781)
    24 = @(908982,330581,794105,854320,1127388,886140,862081,924853,1523979,1182781)
    25 = @(316,115,276,297,391,308,299,321,529,158)
    27 = @(121,44,106,114,150,118,115,123,203,158)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 4 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# Capital row (26): L26/M26 keep their values but switch from the
# "Comma" number style to the same plain style already used by K26
# (and the rest of the row). Copy K26's format onto L26:M26.
# ---------------------------------------------------------------------
$ws.Range("K26").Copy() | Out-Null
$ws.Range("L26:M26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "Workbook updated."
